# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets: F2 305->306, F3 235->236, F5 271->272.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 306
    $ws.Range("F3").Value = 236
    $ws.Range("F5").Value = 272
}
